$d = $word.ActiveDocument
$paraCount = $d.Paragraphs.Count
$eqArrPara = $d.Paragraphs.Item($paraCount - 1)
$insertPos = $eqArrPara.Range.End - 1
$r = $d.Range($insertPos, $insertPos)

$xmlFragment = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid wp14">
        <w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
          <w:sz w:val="40"/>
          <w:szCs w:val="44"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
          <w:sz w:val="40"/>
          <w:szCs w:val="44"/>
        </w:rPr>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
          <w:sz w:val="40"/>
          <w:szCs w:val="44"/>
        </w:rPr>
      </w:pPr>
      <m:oMathPara>
        <m:oMath>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
              <w:sz w:val="40"/>
              <w:szCs w:val="44"/>
            </w:rPr>
            <w:lastRenderedPageBreak/>
            <m:t xml:space="preserve">SD= </m:t>
          </m:r>
          <m:rad>
            <m:radPr>
              <m:degHide m:val="1"/>
              <m:ctrlPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:i/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
              </m:ctrlPr>
            </m:radPr>
            <m:deg/>
            <m:e>
              <m:f>
                <m:fPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:fPr>
                <m:num>
                  <m:nary>
                    <m:naryPr>
                      <m:chr m:val="∑"/>
                      <m:limLoc m:val="undOvr"/>
                      <m:ctrlPr>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                          <w:i/>
                          <w:sz w:val="40"/>
                          <w:szCs w:val="44"/>
                        </w:rPr>
                      </m:ctrlPr>
                    </m:naryPr>
                    <m:sub>
                      <m:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                          <w:sz w:val="40"/>
                          <w:szCs w:val="44"/>
                        </w:rPr>
                        <m:t>i=1</m:t>
                      </m:r>
                    </m:sub>
                    <m:sup>
                      <m:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                          <w:sz w:val="40"/>
                          <w:szCs w:val="44"/>
                        </w:rPr>
                        <m:t>n</m:t>
                      </m:r>
                    </m:sup>
                    <m:e>
                      <m:sSup>
                        <m:sSupPr>
                          <m:ctrlPr>
                            <w:rPr>
                              <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                              <w:i/>
                              <w:sz w:val="40"/>
                              <w:szCs w:val="44"/>
                            </w:rPr>
                          </m:ctrlPr>
                        </m:sSupPr>
                        <m:e>
                          <m:d>
                            <m:dPr>
                              <m:begChr m:val="|"/>
                              <m:endChr m:val="|"/>
                              <m:ctrlPr>
                                <w:rPr>
                                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                                  <w:i/>
                                  <w:sz w:val="40"/>
                                  <w:szCs w:val="44"/>
                                </w:rPr>
                              </m:ctrlPr>
                            </m:dPr>
                            <m:e>
                              <m:sSub>
                                <m:sSubPr>
                                  <m:ctrlPr>
                                    <w:rPr>
                                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                                      <w:i/>
                                      <w:sz w:val="40"/>
                                      <w:szCs w:val="44"/>
                                    </w:rPr>
                                  </m:ctrlPr>
                                </m:sSubPr>
                                <m:e>
                                  <m:r>
                                    <w:rPr>
                                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                                      <w:sz w:val="40"/>
                                      <w:szCs w:val="44"/>
                                    </w:rPr>
                                    <m:t>y</m:t>
                                  </m:r>
                                </m:e>
                                <m:sub>
                                  <m:r>
                                    <w:rPr>
                                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                                      <w:sz w:val="40"/>
                                      <w:szCs w:val="44"/>
                                    </w:rPr>
                                    <m:t>i</m:t>
                                  </m:r>
                                </m:sub>
                              </m:sSub>
                              <m:r>
                                <w:rPr>
                                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                                  <w:sz w:val="40"/>
                                  <w:szCs w:val="44"/>
                                </w:rPr>
                                <m:t>-</m:t>
                              </m:r>
                              <m:acc>
                                <m:accPr>
                                  <m:chr m:val="̅"/>
                                  <m:ctrlPr>
                                    <w:rPr>
                                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                                      <w:i/>
                                      <w:sz w:val="40"/>
                                      <w:szCs w:val="44"/>
                                    </w:rPr>
                                  </m:ctrlPr>
                                </m:accPr>
                                <m:e>
                                  <m:r>
                                    <m:rPr>
                                      <m:scr m:val="script"/>
                                    </m:rPr>
                                    <w:rPr>
                                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                                      <w:sz w:val="40"/>
                                      <w:szCs w:val="44"/>
                                    </w:rPr>
                                    <m:t>Y</m:t>
                                  </m:r>
                                </m:e>
                              </m:acc>
                            </m:e>
                          </m:d>
                        </m:e>
                        <m:sup>
                          <m:r>
                            <w:rPr>
                              <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                              <w:sz w:val="40"/>
                              <w:szCs w:val="44"/>
                            </w:rPr>
                            <m:t>2</m:t>
                          </m:r>
                        </m:sup>
                      </m:sSup>
                    </m:e>
                  </m:nary>
                </m:num>
                <m:den>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>n</m:t>
                  </m:r>
                </m:den>
              </m:f>
            </m:e>
          </m:rad>
        </m:oMath>
      </m:oMathPara>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
          <w:sz w:val="40"/>
          <w:szCs w:val="44"/>
        </w:rPr>
      </w:pPr>
      <m:oMathPara>
        <m:oMath>
          <m:acc>
            <m:accPr>
              <m:chr m:val="̅"/>
              <m:ctrlPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:i/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
              </m:ctrlPr>
            </m:accPr>
            <m:e>
              <m:r>
                <m:rPr>
                  <m:scr m:val="script"/>
                </m:rPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>Y</m:t>
              </m:r>
            </m:e>
          </m:acc>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
              <w:sz w:val="40"/>
              <w:szCs w:val="44"/>
            </w:rPr>
            <m:t>=</m:t>
          </m:r>
          <m:f>
            <m:fPr>
              <m:ctrlPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:i/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
              </m:ctrlPr>
            </m:fPr>
            <m:num>
              <m:sSub>
                <m:sSubPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:sSubPr>
                <m:e>
                  <m:r>
                    <m:rPr>
                      <m:scr m:val="script"/>
                    </m:rPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>y</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>1</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>+</m:t>
              </m:r>
              <m:sSub>
                <m:sSubPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:sSubPr>
                <m:e>
                  <m:r>
                    <m:rPr>
                      <m:scr m:val="script"/>
                    </m:rPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>y</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>2</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>+…+</m:t>
              </m:r>
              <m:sSub>
                <m:sSubPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:sSubPr>
                <m:e>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>y</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>N</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
            </m:num>
            <m:den>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>N</m:t>
              </m:r>
            </m:den>
          </m:f>
        </m:oMath>
      </m:oMathPara>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
          <w:sz w:val="40"/>
          <w:szCs w:val="44"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
          <w:sz w:val="40"/>
          <w:szCs w:val="44"/>
        </w:rPr>
      </w:pPr>
      <m:oMathPara>
        <m:oMath>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
              <w:sz w:val="40"/>
              <w:szCs w:val="44"/>
            </w:rPr>
            <m:t>Y=</m:t>
          </m:r>
          <m:d>
            <m:dPr>
              <m:begChr m:val="{"/>
              <m:endChr m:val="}"/>
              <m:ctrlPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:i/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
              </m:ctrlPr>
            </m:dPr>
            <m:e>
              <m:sSub>
                <m:sSubPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:sSubPr>
                <m:e>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>y</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>1</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>,</m:t>
              </m:r>
              <m:sSub>
                <m:sSubPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:sSubPr>
                <m:e>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>y</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>2</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>,</m:t>
              </m:r>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>…</m:t>
              </m:r>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>,</m:t>
              </m:r>
              <m:sSub>
                <m:sSubPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:sSubPr>
                <m:e>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>y</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>n</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
            </m:e>
          </m:d>
        </m:oMath>
      </m:oMathPara>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
          <w:sz w:val="40"/>
          <w:szCs w:val="44"/>
        </w:rPr>
      </w:pPr>
      <m:oMathPara>
        <m:oMath>
          <m:r>
            <m:rPr>
              <m:scr m:val="script"/>
            </m:rPr>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
              <w:sz w:val="40"/>
              <w:szCs w:val="44"/>
            </w:rPr>
            <m:t>Y=</m:t>
          </m:r>
          <m:d>
            <m:dPr>
              <m:begChr m:val="{"/>
              <m:endChr m:val="}"/>
              <m:ctrlPr>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:i/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
              </m:ctrlPr>
            </m:dPr>
            <m:e>
              <m:sSub>
                <m:sSubPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:sSubPr>
                <m:e>
                  <m:r>
                    <m:rPr>
                      <m:scr m:val="script"/>
                    </m:rPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>y</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>1</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>,</m:t>
              </m:r>
              <m:sSub>
                <m:sSubPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:sSubPr>
                <m:e>
                  <m:r>
                    <m:rPr>
                      <m:scr m:val="script"/>
                    </m:rPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>y</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>2</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
              <m:r>
                <w:rPr>
                  <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                  <w:sz w:val="40"/>
                  <w:szCs w:val="44"/>
                </w:rPr>
                <m:t>,…,</m:t>
              </m:r>
              <m:sSub>
                <m:sSubPr>
                  <m:ctrlPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:i/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                  </m:ctrlPr>
                </m:sSubPr>
                <m:e>
                  <m:r>
                    <m:rPr>
                      <m:scr m:val="script"/>
                    </m:rPr>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>y</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <w:rPr>
                      <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
                      <w:sz w:val="40"/>
                      <w:szCs w:val="44"/>
                    </w:rPr>
                    <m:t>N</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
            </m:e>
          </m:d>
        </m:oMath>
      </m:oMathPara>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsiaTheme="minorEastAsia"/>
          <w:sz w:val="40"/>
          <w:szCs w:val="44"/>
        </w:rPr>
      </w:pPr>
      <m:oMathPara>
        <m:oMath>
          <m:r>
            <w:rPr>
              <w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/>
              <w:sz w:val="40"/>
              <w:szCs w:val="44"/>
            </w:rPr>
            <m:t>N=n</m:t>
          </m:r>
        </m:oMath>
      </m:oMathPara>
    </w:p>

        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xmlFragment)
Write-Output "Inserted. New paragraph count:"
Write-Output $d.Paragraphs.Count
